$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the "Puerro" (leek) sheet.
# It belongs chronologically where the old row 70 used to sit, so insert a
# fresh row there and push every following record down by one.
$ws.Rows("70").Insert()

$ws.Range("A70").Value = 10
$ws.Range("B70").Value = "Vega Modelo de Temuco"
$ws.Range("C70").Value = "La Araucanía"
$ws.Range("D70").Value = 44662
$ws.Range("E70").Value = 9
$ws.Range("F70").Value = 100112005
$ws.Range("G70").Value = "Puerro"
$ws.Range("H70").Value = "Azul de Maquehue"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 40
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = 12000
$ws.Range("N70").Value = "$/docena de paquetes"
$ws.Range("O70").Value = "Provincia de Cautín"
$ws.Range("P70").Value = 1000
$ws.Range("Q70").Value = 12
$ws.Range("R70").Value = "Hortaliza"
